$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the hyperlinks on A2/A3 (TFVC sample URLs) entirely.
$ws.Hyperlinks.Delete()

# Row 2 (sample "project1" row) is removed outright -- row 3 keeps its row
# number (it is NOT shifted up to become row 2). Achieve this by deleting
# row 2 (which shifts row 3 up to row 2), then re-inserting a blank row at
# position 2 (which pushes the old row 3's content back down to row 3),
# and finally clearing that freshly inserted blank row completely so it
# leaves no trace (no cells at all) once saved.
$ws.Rows(2).Delete()
$ws.Rows(2).Insert()
$ws.Range("A2:C2").Clear()

# Row 3 (the other sample "project2" row) keeps its styling but loses its
# sample values -- it becomes a blank input row.
$ws.Range("A3:C3").ClearContents()

# Add the two new header columns, copying the header style (bold white on
# dark-blue fill, centered + wrap) from the existing PAT header cell so the
# same cell style gets reused.
$ws.Range("D1").Value = "Username"
$ws.Range("E1").Value = "Password"
$ws.Range("C1").Copy()
$ws.Range("D1:E1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Restore the saved cursor/selection position recorded in the new file.
$ws.Range("C12").Select()
